# Update symbol list (cryptos price snapshot) with refreshed values.
# The "D" column holds prices stored as text (inline strings that look like
# numbers, e.g. "245.03", "0.03270"), so we force the cell format to Text
# ("@") before assigning the new value. This prevents Excel from converting
# the string into a floating point number and losing significant trailing
# zeros (e.g. "3.500" -> 3.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    "D2"  = "245.01"
    "D4"  = "5.007"
    "D5"  = "0.05613"
    "D7"  = "3.005"
    "D8"  = "0.8107"
    "D9"  = "0.8402"
    "D10" = "0.1337"
    "D11" = "0.06949"
    "D13" = "0.09409"
    "D14" = "0.001516"
    "D15" = "0.0005967"
    "D16" = "0.006085"
    "D17" = "3.500"
    "D18" = "2.092"
    "D20" = "0.03282"
    "D22" = "3.740"
    "D23" = "0.04673"
    "D24" = "0.1369"
    "D26" = "0.004523"
    "D27" = "0.00009698"
    "D28" = "0.0001940"
    "D41" = "0.1359"
    "D42" = "0.006239"
    "D43" = "0.002735"
    "D44" = "0.008060"
    "D45" = "0.00005273"
    "D47" = "0.1800"
    "D48" = "0.002042"
    "D49" = "0.00002100"
    "D50" = "0.0002000"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) column text updates (plain text, not numeric so no special
# formatting is required to keep them as strings).
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("E42").Value = "41KickTokenKICKBestin24h"
